$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update courier ID number for existing row (row 2)
$ws.Range("B2").Value = 123321

# Copy formatting from row 2 to new row 3 so styles match
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)

# Populate new row 3 with the new courier's data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 5363124
$ws.Range("C3").Value = "Neman Ismiyev"
$ws.Range("D3").Value = 134.15
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = "-"

# Match final selection cell from the saved workbook
$ws.Range("D8").Select()
